# UPDATE technology portfolios for Norway
# Set the base (2025) total cost value in F2, the other sheets (2030-2050)
# recompute automatically since they hold formulas referencing '2025'!F2.

$wb = $excel.ActiveWorkbook
$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("F2").Value = 1100000
